$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.494.98'
$ws.Range('E2').Value = '  -0.96%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.582.03'
$ws.Range('E3').Value = '  -1.90%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.80'
$ws.Range('E5').Value = '  -1.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.20'
$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('E8').Value = '  -1.71%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.581.29'
$ws.Range('E9').Value = '  -1.91%  '

$ws.Range('E10').Value = '  -4.14%  '

$ws.Range('E11').Value = '  +0.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  -1.53%  '

$ws.Range('E13').Value = '  -1.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.69'
$ws.Range('E14').Value = '  -3.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.052.24'
$ws.Range('E15').Value = '  -1.89%  '

$ws.Range('E16').Value = '  -2.32%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.357.91'
$ws.Range('E17').Value = '  -0.97%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.595.76'
$ws.Range('E18').Value = '  -0.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.42'
$ws.Range('E19').Value = '  -5.79%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.72'
$ws.Range('E20').Value = '  -3.96%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.44'
$ws.Range('E21').Value = '  -2.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.22'
$ws.Range('E22').Value = '  -2.98%  '

$ws.Range('E23').Value = '  -1.60%  '

$ws.Range('E24').Value = '  +0.02%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.88'
$ws.Range('E25').Value = '  -3.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.02'
$ws.Range('E26').Value = '  -8.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '68.82'
$ws.Range('E27').Value = '  -2.87%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.713.54'
$ws.Range('E28').Value = '  -1.91%  '

$ws.Range('E29').Value = '  +0.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0985'
$ws.Range('E30').Value = '  -2.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '534.99'
$ws.Range('E31').Value = '  -3.49%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.00'
$ws.Range('E32').Value = '  +0.59%  '

$ws.Range('E33').Value = '  -2.92%  '

$ws.Range('E34').Value = '  -2.43%  '

$ws.Range('E35').Value = '  -1.80%  '

$ws.Range('E36').Value = '  +0.04%  '

$ws.Range('E37').Value = '  -3.45%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.64'
$ws.Range('E38').Value = '  -0.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.73'
$ws.Range('E39').Value = '  -2.43%  '

$ws.Range('E40').Value = '  -1.97%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.25'
$ws.Range('E41').Value = '  +1.92%  '

$ws.Range('E42').Value = '  -0.59%  '

$ws.Range('E43').Value = '  -2.09%  '

$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('E45').Value = '  -2.46%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0285'
$ws.Range('E46').Value = '  -4.77%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.08'
$ws.Range('E47').Value = '  -2.04%  '

$ws.Range('E48').Value = '  -3.54%  '

$ws.Range('E49').Value = '  -2.51%  '

$ws.Range('E50').Value = '  -1.06%  '

$ws.Range('E51').Value = '  -1.63%  '
